# Week 17 data log + tiebreaking fix (commit: "Logged Week 17 data and
# fixed Simulate_Season.py tiebreaking method").
#
# Appends this week's per-play yardage samples to the YDS / ST running
# logs, and rolls the new week's counts into the OFF / DEF / ST / TURNS /
# PEN season totals.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# YDS sheet - append this week's individual play-yardage samples to the
# running space-separated lists.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("YDS")

$ws.Cells.Item(2,2).Value2 = $ws.Cells.Item(2,2).Value2 + " 4 -1 2 0 -1 5 4 6 3 3 1 0 12 -1 4 1"
$ws.Cells.Item(2,3).Value2 = $ws.Cells.Item(2,3).Value2 + " 24 14 1 2 1 2 2 1 -5 2 7 4 1 -3 2 7 12 5 0 1 2 2 1 1 15 3 11 6 9"
$ws.Cells.Item(3,2).Value2 = $ws.Cells.Item(3,2).Value2 + " 9 6 7 7 14 11 9 6 9 7 7 11 21 12 5 26 0 8 31 16 8 4"
$ws.Cells.Item(3,3).Value2 = $ws.Cells.Item(3,3).Value2 + " 42 4 17 3 3 3 20 10 5 3 4 9 4 5 5 11 43 -7 6 19 16 7 12 7 11"

# ---------------------------------------------------------------------
# OFF sheet - season offensive totals through Week 17.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("OFF")

$ws.Cells.Item(2,3).Value2  = 211
$ws.Cells.Item(2,5).Value2  = 4
$ws.Cells.Item(2,6).Value2  = 59
$ws.Cells.Item(2,7).Value2  = 54
$ws.Cells.Item(2,8).Value2  = 5
$ws.Cells.Item(2,10).Value2 = 32
$ws.Cells.Item(2,12).Value2 = 280
$ws.Cells.Item(2,13).Value2 = 188
$ws.Cells.Item(2,15).Value2 = 22
$ws.Cells.Item(2,16).Value2 = 12
$ws.Cells.Item(2,17).Value2 = 517

$ws.Cells.Item(3,2).Value2  = 14
$ws.Cells.Item(3,3).Value2  = 193
$ws.Cells.Item(3,5).Value2  = 36
$ws.Cells.Item(3,6).Value2  = 135
$ws.Cells.Item(3,7).Value2  = 45
$ws.Cells.Item(3,8).Value2  = 22
$ws.Cells.Item(3,9).Value2  = 58
$ws.Cells.Item(3,10).Value2 = 62
$ws.Cells.Item(3,14).Value2 = 13

# ---------------------------------------------------------------------
# DEF sheet - season defensive totals through Week 17.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("DEF")

$ws.Cells.Item(2,2).Value2  = 5
$ws.Cells.Item(2,3).Value2  = 191
$ws.Cells.Item(2,4).Value2  = 10
$ws.Cells.Item(2,6).Value2  = 59
$ws.Cells.Item(2,7).Value2  = 38
$ws.Cells.Item(2,9).Value2  = 6
$ws.Cells.Item(2,10).Value2 = 25
$ws.Cells.Item(2,12).Value2 = 279
$ws.Cells.Item(2,13).Value2 = 166
$ws.Cells.Item(2,15).Value2 = 22
$ws.Cells.Item(2,16).Value2 = 13
$ws.Cells.Item(2,17).Value2 = 507

$ws.Cells.Item(3,3).Value2  = 188
$ws.Cells.Item(3,5).Value2  = 39
$ws.Cells.Item(3,6).Value2  = 114
$ws.Cells.Item(3,7).Value2  = 36
$ws.Cells.Item(3,9).Value2  = 68
$ws.Cells.Item(3,10).Value2 = 54
$ws.Cells.Item(3,14).Value2 = 22

# ---------------------------------------------------------------------
# ST sheet - special-teams season totals + appended per-kick logs.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ST")

$ws.Cells.Item(2,2).Value2  = 97
$ws.Cells.Item(2,4).Value2  = 62
$ws.Cells.Item(2,6).Value2  = 352
$ws.Cells.Item(2,7).Value2  = 338
$ws.Cells.Item(2,8).Value2  = 10
$ws.Cells.Item(2,9).Value2  = 6
$ws.Cells.Item(2,12).Value2 = 81

$ws.Cells.Item(3,2).Value2  = 82

$ws.Cells.Item(3,4).Value2 = $ws.Cells.Item(3,4).Value2 + " 42 57 45 51"
$ws.Cells.Item(4,4).Value2 = $ws.Cells.Item(4,4).Value2 + " 0 -5 0 0"
$ws.Cells.Item(5,4).Value2 = $ws.Cells.Item(5,4).Value2 + " 0 17 -2"
$ws.Cells.Item(6,2).Value2 = $ws.Cells.Item(6,2).Value2 + " 20 18 20"

# ---------------------------------------------------------------------
# TURNS sheet - season turnover totals.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("TURNS")

$ws.Cells.Item(2,2).Value2 = 3
$ws.Cells.Item(2,4).Value2 = 10
$ws.Cells.Item(3,4).Value2 = 8
$ws.Cells.Item(3,5).Value2 = 4

# ---------------------------------------------------------------------
# PEN sheet - season penalty totals.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("PEN")

$ws.Cells.Item(2,2).Value2 = 17
$ws.Cells.Item(2,4).Value2 = 15
$ws.Cells.Item(3,2).Value2 = 29
$ws.Cells.Item(4,4).Value2 = 10
